$wb = $excel.ActiveWorkbook

# --- Fix the cell values on "VerifyStartEndDateValidation" ---
$ws6 = $wb.Worksheets.Item("VerifyStartEndDateValidation")

# I2: typo fix "End Dat ... Start Dat!" -> "End Date ... Start Date!"
$ws6.Range("I2").Value = "End Date should always be greater or equal to the Start Date!"

# G2: was "24/03/2017" -> "24/04/2017"
# (leading apostrophe keeps it stored/styled as text, matching the original
#  quote-prefixed date-formatted cell instead of letting it re-parse as a date)
$ws6.Range("G2").Value = "'24/04/2017"

# H2: was "Test User" -> "test demo"
# (leading apostrophe preserves the original quote-prefixed text style)
$ws6.Range("H2").Value = "'test demo"

# --- Switch the active/selected tab from "Test Cases" to "VerifyStartEndDateValidation" ---
[void]$ws6.Activate()
[void]$ws6.Range("D2").Select()
